# ---------------------------------------------------------------------------
# Applies "Jun's file updates for all IO data and others" to the workbook:
#   - About sheet: updated source citation (year/title/URL/page ref)
#   - EIA Table 1 sheet: refreshed elasticity figures (AEO2003 source data)
#   - EoCEDwEC sheet: formulas now diff against the 1-year column (B) instead
#     of the 3-year column (D), and the active tab / selections reset back to
#     the "About" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsEia   = $wb.Worksheets.Item("EIA Table 1")
$wsEoC   = $wb.Worksheets.Item("EoCEDwEC")

# --- About sheet: new source citation -------------------------------------
$wsAbout.Range("B4").Value = 2005
$wsAbout.Range("B5").Value = "Price Responsiveness in the AEO 2003 NEMS Residential and Commercial Buildings Sector Models"
$wsAbout.Range("B6").Value = "http://www.eia.gov/oiaf/analysispaper/elasticity/pdf/tbl.pdf"
$wsAbout.Range("B7").Value = "Page 1, Table 1"

$wsAbout.Rows.Item(4).RowHeight = 14.45
$wsAbout.Rows.Item(5).RowHeight = 14.45
$wsAbout.Rows.Item(6).RowHeight = 14.45
$wsAbout.Rows.Item(7).RowHeight = 14.45

# --- EIA Table 1 sheet: refreshed elasticity figures -----------------------
# Urban Residential block (rows 7-9)
$wsEia.Range("B7").Value = -0.2
$wsEia.Range("C7").Value = -0.28999999999999998
$wsEia.Range("D7").Value = -0.34
$wsEia.Range("E7").Value = -0.49
$wsEia.Range("F7").Value = 0.01
$wsEia.Range("G7").Value = 0

$wsEia.Range("B8").Value = -0.14000000000000001
$wsEia.Range("C8").Value = -0.24
$wsEia.Range("D8").Value = -0.3
$wsEia.Range("E8").Value = 0.13
$wsEia.Range("F8").Value = -0.41
$wsEia.Range("G8").Value = 0.02

$wsEia.Range("B9").Value = -0.15
$wsEia.Range("C9").Value = -0.27
$wsEia.Range("D9").Value = -0.34
$wsEia.Range("E9").Value = 0.01
$wsEia.Range("F9").Value = 0.05
$wsEia.Range("G9").Value = -0.6

# Rural Residential block (rows 14-16)
$wsEia.Range("B14").Value = -0.1
$wsEia.Range("C14").Value = -0.17
$wsEia.Range("D14").Value = -0.2
$wsEia.Range("E14").Value = -0.45
$wsEia.Range("F14").Value = 0.01
$wsEia.Range("G14").Value = 0

$wsEia.Range("B15").Value = -0.14000000000000001
$wsEia.Range("C15").Value = -0.24
$wsEia.Range("D15").Value = -0.28999999999999998
$wsEia.Range("E15").Value = 0.86
$wsEia.Range("F15").Value = -0.4
$wsEia.Range("G15").Value = 0.01

$wsEia.Range("B16").Value = -0.13
$wsEia.Range("C16").Value = -0.23
$wsEia.Range("D16").Value = -0.28000000000000003
$wsEia.Range("E16").Value = 0.08
$wsEia.Range("F16").Value = 0.75
$wsEia.Range("G16").Value = -0.39

# --- EoCEDwEC sheet: formulas now reference column B instead of column D ---
$wsEoC.Range("B2").Formula = "='EIA Table 1'!E7-'EIA Table 1'!B7"
$wsEoC.Range("D2").Formula = "='EIA Table 1'!E14-'EIA Table 1'!B14"
$wsEoC.Range("B4").Formula = "='EIA Table 1'!F8-'EIA Table 1'!B8"
$wsEoC.Range("D4").Formula = "='EIA Table 1'!F15-'EIA Table 1'!B15"
$wsEoC.Range("B5").Formula = "='EIA Table 1'!G9-'EIA Table 1'!B9"
$wsEoC.Range("D5").Formula = "='EIA Table 1'!G16-'EIA Table 1'!B16"

# --- Active tab / selection: back to "About" --------------------------------
$wsEia.Activate()
$wsEia.Range("A1").Select()

$wsEoC.Activate()
$wsEoC.Range("A1").Select()

$wsAbout.Activate()
$wsAbout.Range("A1").Select()
